# Updated symbol list on Tue Dec 13 16:54:43 UTC 2022 with GitHub Actions
#
# Applies the "Price" (column D) refresh for most rows, plus a full
# row-content swap between row 42 (CEJI -> BKEXToken) and row 43
# (BKEXToken -> CEJI), matching the upstream coinranking.com scrape diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A reference cell whose style (xf index 0 / "General", no quote-prefix)
# we reuse after writing numeric-looking text, so the cell's style index
# is not perturbed by forcing a text number format.
$styleDonor = $ws.Range("A1")

function Set-TextValue {
    # Positional params only - named args (-Cell/-Text) are not bound by
    # this PowerShell-style interpreter.
    param($Cell, $Text)
    # Force the literal text (even if it looks like a number, e.g. "23.01")
    # to be stored as a text value rather than being auto-coerced to a
    # number by Excel's input parser.
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    # Restore the original (default) cell style so the style index isn't
    # left pointing at the throwaway text-format xf record we just used.
    $Cell.Style = $styleDonor.Style
}

# --- Column D ("Price") refreshes -----------------------------------
$priceUpdates = @{
    3  = "23.01"
    4  = "6.318"
    5  = "0.06227"
    6  = "3.649"
    7  = "6.674"
    8  = "1.400"
    9  = "0.8332"
    10 = "0.01383"
    11 = "0.1607"
    12 = "0.08316"
    13 = "0.03559"
    14 = "0.03187"
    15 = "4.058"
    16 = "0.09311"
    17 = "0.001641"
    18 = "0.04742"
    19 = "0.006356"
    20 = "0.005699"
    22 = "0.0001501"
    24 = "2.325"
    25 = "0.3357"
    27 = "0.0002706"
    40 = "0.04737"
    41 = "0.006981"
    44 = "0.01215"
    45 = "0.00006275"
    46 = "0.0009899"
    48 = "0.7826"
    49 = "0.002348"
    50 = "0.00002402"
    51 = "0.01241"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    Set-TextValue $cell $priceUpdates[$row]
}

# --- Rows 42 / 43: CEJI and BKEXToken swap places --------------------
# Row 42 becomes BKEXToken, row 43 becomes CEJI (with refreshed prices).
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D42") "0.1168"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D43") "0.003802"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"
